$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to make edits, then restore protection afterwards.
$ws.Unprotect()

# Update the disclosure text date from 2021-04-28 to 2021-04-29
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1044284928778688
$ws.Range("E2").Value = -0.02118366455557996

$ws.Range("D3").Value = 0.1092060203803956
$ws.Range("E3").Value = -0.004359047892170786

$ws.Range("D4").Value = 0.1166100878626287
$ws.Range("E4").Value = 0.009342125922436884

$ws.Range("D5").Value = 0.1379779879567954
$ws.Range("E5").Value = 0.003744657274022511

$ws.Range("D6").Value = 0.1324825545163154
$ws.Range("E6").Value = 0.005018820577164407

$ws.Range("D7").Value = 0.14101985338517
$ws.Range("E7").Value = 0.006409643164641876

$ws.Range("D8").Value = 0.1292873109736802
$ws.Range("E8").Value = 0.005428571428571338

$ws.Range("D9").Value = 0.1289876920471459
$ws.Range("E9").Value = 0.0007404239818180613

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.00128399805780588

# Restore sheet protection (content-protected, matching the original intent).
$ws.Protect()
